$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.465.13'
$ws.Range("E2").Value = '  -0.89%  '
$ws.Range("D3").Value = '3.472.42'
$ws.Range("E3").Value = '  -1.14%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '588.73'
$ws.Range("E5").Value = '  +2.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.80'
$ws.Range("E6").Value = '  -1.67%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.605'
$ws.Range("E7").Value = '  -1.87%  '
$ws.Range("D8").Value = '3.467.24'
$ws.Range("E8").Value = '  -1.15%  '
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.190'
$ws.Range("E10").Value = '  +0.80%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.76'
$ws.Range("E11").Value = '  +1.90%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.571'
$ws.Range("E12").Value = '  -4.32%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '46.43'
$ws.Range("E13").Value = '  -1.34%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000276'
$ws.Range("E14").Value = '  +0.53%  '
$ws.Range("D15").Value = '4.033.44'
$ws.Range("E15").Value = '  -0.97%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '613.48'
$ws.Range("E16").Value = '  -10.26%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '8.29'
$ws.Range("E17").Value = '  -4.85%  '
$ws.Range("D18").Value = '3.492.65'
$ws.Range("E18").Value = '  -0.25%  '
$ws.Range("D19").Value = '68.476.53'
$ws.Range("E19").Value = '  -0.85%  '
$ws.Range("E20").Value = '  -2.18%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.18'
$ws.Range("E21").Value = '  -1.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.08'
$ws.Range("E22").Value = '  -0.28%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.869'
$ws.Range("E23").Value = '  -4.33%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '15.77'
$ws.Range("E24").Value = '  -4.60%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '95.54'
$ws.Range("E25").Value = '  -1.96%  '
$ws.Range("E26").Value = '  -1.52%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  -0.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.60'
$ws.Range("E28").Value = '  -2.23%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.08'
$ws.Range("E29").Value = '  -3.61%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.67'
$ws.Range("E30").Value = '  -1.35%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.38'
$ws.Range("E31").Value = '  -4.93%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.07'
$ws.Range("E32").Value = '  -3.08%  '
$ws.Range("E33").Value = '  -2.41%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.79'
$ws.Range("E34").Value = '  -6.66%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '571.21'
$ws.Range("E35").Value = '  -1.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.65'
$ws.Range("E36").Value = '  -1.60%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.48'
$ws.Range("E37").Value = '  -5.71%  '
$ws.Range("E38").Value = '  -0.92%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.101'
$ws.Range("E39").Value = '  -3.88%  '
$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  +0.08%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.136'
$ws.Range("E41").Value = '  -0.76%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0435'
$ws.Range("E42").Value = '  -0.79%  '
$ws.Range("D43").Value = '3.388.76'
$ws.Range("E43").Value = '  -2.15%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.322'
$ws.Range("E44").Value = '  -4.36%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '32.53'
$ws.Range("E45").Value = '  -2.10%  '
$ws.Range("D46").Value = '0.0₃0689'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.80'
$ws.Range("E47").Value = '  -2.87%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.54'
$ws.Range("E48").Value = '  -1.54%  '
$ws.Range("E49").Value = '  -3.18%  '
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '132.46'
$ws.Range("E50").Value = '  -0.94%  '
$ws.Range("B51").Value = 'MXToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.67'
$ws.Range("E51").Value = '  +12.25%  '
